# Append 5 new service-log rows (284-288) to Sheet1, matching existing
# table layout: A=DATE, B=VECHILE REG NO, C=VEHICLE BRAND, D=ISSUE,
# E=STATUS, F=AMOUNT, G=CASH TYPE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 284
$ws.Cells.Item(284,1).Value = 44792
$ws.Cells.Item(284,2).Value = "KA03MQ5430"
$ws.Cells.Item(284,3).Value = "I20 ASTA"
$ws.Cells.Item(284,4).Value = "GENERAL CHECKUP"
$ws.Cells.Item(284,5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(284,6).Value = 1487
$ws.Cells.Item(284,7).Value = "CREDIT"

# Row 285
$ws.Cells.Item(285,1).Value = 44792
$ws.Cells.Item(285,2).Value = "KA03MZ9550"
$ws.Cells.Item(285,3).Value = "ECOSPORT"
$ws.Cells.Item(285,4).Value = "MIRROR CHANGE"
$ws.Cells.Item(285,5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(285,6).Value = 5416
$ws.Cells.Item(285,7).Value = "CREDIT"

# Row 286
$ws.Cells.Item(286,1).Value = 44792
$ws.Cells.Item(286,2).Value = "KA03MV0364"
$ws.Cells.Item(286,3).Value = "SCORPIO"
$ws.Cells.Item(286,4).Value = "PMS                                      WW"
$ws.Cells.Item(286,5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(286,6).Value = 12737
$ws.Cells.Item(286,7).Value = "CREDIT"

# Row 287
$ws.Cells.Item(287,1).Value = 44792
$ws.Cells.Item(287,2).Value = "KA03MQ5430"
$ws.Cells.Item(287,3).Value = "I20 ASTA"
$ws.Cells.Item(287,4).Value = "GENERAL CHECKUP         WW"
$ws.Cells.Item(287,5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(287,6).Value = 1487
$ws.Cells.Item(287,7).Value = "CREDIT"

# Row 288
$ws.Cells.Item(288,1).Value = 44792
$ws.Cells.Item(288,2).Value = "WB20Z5652"
$ws.Cells.Item(288,3).Value = "FIGO"
$ws.Cells.Item(288,4).Value = "GENERAL CHECKUP"
$ws.Cells.Item(288,5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(288,6).Value = 1010
$ws.Cells.Item(288,7).Value = "P PAY"

# Select the last entered cell, matching the author's final cursor position
$ws.Range("G288").Select()
